$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add C1 date (13-01-2023), copying B1 bold/border/center style ---
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("C1").Value = "13-01-2023"

# --- Data rows 2-36: labels re-sorted alphabetically, new column C values added ---
$ws.Range("A2").Value = "1810 Renta variable"
$ws.Range("B2").Value = 4113761.49
$ws.Range("C2").Value = 4114314.89
$ws.Range("A3").Value = "1822 Raices Valores Negociables"
$ws.Range("B3").Value = 3706367.68
$ws.Range("C3").Value = 3723756.54
$ws.Range("A4").Value = "Allaria Acciones"
$ws.Range("B4").Value = 474972.72
$ws.Range("C4").Value = 472594.6
$ws.Range("A5").Value = "Alpha Acciones"
$ws.Range("B5").Value = 2243785.92
$ws.Range("C5").Value = 2249629.99
$ws.Range("A6").Value = "Alpha Mega"
$ws.Range("B6").Value = 4547640.7
$ws.Range("C6").Value = 4562501.4
$ws.Range("A7").Value = "Argenfunds"
$ws.Range("B7").Value = 143994.05
$ws.Range("C7").Value = 143868.09
$ws.Range("A8").Value = "Balanz"
$ws.Range("B8").Value = 488559.75
$ws.Range("C8").Value = 482158.04
$ws.Range("A9").Value = "Delta Acciones"
$ws.Range("B9").Value = 280180.48
$ws.Range("C9").Value = 282838.94
$ws.Range("A10").Value = "Delta Internacional"
$ws.Range("B10").Value = 156779.61
$ws.Range("C10").Value = 156583.47
$ws.Range("A11").Value = "Delta Latinoamerica"
$ws.Range("B11").Value = 250911.45
$ws.Range("C11").Value = 250949.93
$ws.Range("A12").Value = "Delta Recursos Naturales"
$ws.Range("B12").Value = 5591096.83
$ws.Range("C12").Value = 5590688.21
$ws.Range("A13").Value = "Delta Select"
$ws.Range("B13").Value = 4893135.39
$ws.Range("C13").Value = 4897975.6
$ws.Range("A14").Value = "FBA Acciones Argentinas"
$ws.Range("B14").Value = 986827.87
$ws.Range("C14").Value = 984282.63
$ws.Range("A15").Value = "FBA Calificado"
$ws.Range("B15").Value = 958088.42
$ws.Range("C15").Value = 964510.55
$ws.Range("A16").Value = "Fima Acciones"
$ws.Range("B16").Value = 1238457.64
$ws.Range("C16").Value = 1232401.18
$ws.Range("A17").Value = "Fima PB Acciones"
$ws.Range("B17").Value = 2689607.22
$ws.Range("C17").Value = 2746566.57
$ws.Range("A18").Value = "Galileo Acciones"
$ws.Range("B18").Value = 34978162.49
$ws.Range("C18").Value = 35247486.51
$ws.Range("A19").Value = "Goal Acciones Argentinas"
$ws.Range("B19").Value = 571901.09
$ws.Range("C19").Value = 572507.54
$ws.Range("A20").Value = "Goal acciones plus"
$ws.Range("B20").Value = 10033.82
$ws.Range("C20").Value = 10218.84
$ws.Range("A21").Value = "HF Acciones Argentinas"
$ws.Range("B21").Value = 1110944.81
$ws.Range("C21").Value = 1114191.59
$ws.Range("A22").Value = "HF Acciones Lideres"
$ws.Range("B22").Value = 2532721.25
$ws.Range("C22").Value = 2532858.42
$ws.Range("A23").Value = "IAM Renta Variable"
$ws.Range("B23").Value = 467361.57
$ws.Range("C23").Value = 488019.84
$ws.Range("A24").Value = "IEB Value"
$ws.Range("B24").Value = 53505.18
$ws.Range("C24").Value = 53562.38
$ws.Range("A25").Value = "Lombardi"
$ws.Range("B25").Value = 527821.1800000001
$ws.Range("C25").Value = 577488.4300000001
$ws.Range("A26").Value = "Megainver"
$ws.Range("B26").Value = 404728.96
$ws.Range("C26").Value = 405310.49
$ws.Range("A27").Value = "Pellegrini Acciones"
$ws.Range("B27").Value = 936650.5600000001
$ws.Range("C27").Value = 940287.5
$ws.Range("A28").Value = "Pionero Acciones"
$ws.Range("B28").Value = 1411757.85
$ws.Range("C28").Value = 1407218.71
$ws.Range("A29").Value = "Premier Renta Variable"
$ws.Range("B29").Value = 29979.9
$ws.Range("C29").Value = 29113.3
$ws.Range("A30").Value = "Quinquela Acciones"
$ws.Range("B30").Value = 1214413.21
$ws.Range("C30").Value = 1216155.65
$ws.Range("A31").Value = "Rofex 20 Renta Variable"
$ws.Range("B31").Value = 925134.17
$ws.Range("C31").Value = 925542.4399999999
$ws.Range("A32").Value = "Supefondo RV"
$ws.Range("B32").Value = 17717880.06
$ws.Range("C32").Value = 17725718.49
$ws.Range("A33").Value = "Superfondo "
$ws.Range("B33").Value = 24944544.43
$ws.Range("C33").Value = 24951084.52
$ws.Range("A34").Value = "Supergestion"
$ws.Range("B34").Value = 4062543.59
$ws.Range("C34").Value = 4051087.37
$ws.Range("A35").Value = "avg"
$ws.Range("B35").Value = 3777704.59
$ws.Range("C35").Value = 3791014.32
$ws.Range("A36").Value = "total"
$ws.Range("B36").Value = 124664251.34
$ws.Range("C36").Value = 125103472.65
